# Applies the "Completed more 10K Steps trials" edit:
#  - Row 20 (Epoch 1 Min Loss) gets filled in with real values
#  - Rows 21-28 become "Loss 1".."Loss 8" rows with new data
#  - The old "Epoch 2 Min Loss" / "Output Grade (A-F)" rows move down to rows
#    29/30 and the "Epoch 2 Min Loss" row gets a value filled in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlNone = -4142
$xlThin = 2
$xlMedium = -4138

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# A pristine, never-touched cell we use as a source of "no border / default
# style" formatting so that re-using a previously styled cell never drags
# along stray border colors.
$blank = "Z100"

# ---------------------------------------------------------------------------
# Build each brand-new border combination once on a scratch cell far away
# from the used range, then stamp it wherever it's needed with
# PasteSpecial(xlPasteFormats). Building on a pristine cell (instead of
# editing an already-bordered cell in place) avoids the engine dragging
# along stray colors from the cell's previous border.
# ---------------------------------------------------------------------------

# Style A: top = medium only  (used by E20)
$styleA = "Z101"
$ws.Range($styleA).Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range($styleA).Borders.Item($xlEdgeTop).Weight = $xlMedium

# Style B: left = thin, right = thin, top = medium (used by F20:I20)
$styleB = "Z102"
$ws.Range($styleB).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range($styleB).Borders.Item($xlEdgeLeft).Weight = $xlThin
$ws.Range($styleB).Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range($styleB).Borders.Item($xlEdgeRight).Weight = $xlThin
$ws.Range($styleB).Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range($styleB).Borders.Item($xlEdgeTop).Weight = $xlMedium

# Style C: left = thin, right = thin only (used by F21:I28 and F29:I29)
$styleC = "Z103"
$ws.Range($styleC).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range($styleC).Borders.Item($xlEdgeLeft).Weight = $xlThin
$ws.Range($styleC).Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range($styleC).Borders.Item($xlEdgeRight).Weight = $xlThin

# Style D: left = medium only (used by E29)
$styleD = "Z104"
$ws.Range($styleD).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range($styleD).Borders.Item($xlEdgeLeft).Weight = $xlMedium

# Style E: left = medium, bottom = thin (used by E30)
$styleE = "Z105"
$ws.Range($styleE).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range($styleE).Borders.Item($xlEdgeLeft).Weight = $xlMedium
$ws.Range($styleE).Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range($styleE).Borders.Item($xlEdgeBottom).Weight = $xlThin

# Style F: left = thin, right = thin, bottom = thin (used by F30:I30)
$styleF = "Z106"
$ws.Range($styleF).Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range($styleF).Borders.Item($xlEdgeLeft).Weight = $xlThin
$ws.Range($styleF).Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range($styleF).Borders.Item($xlEdgeRight).Weight = $xlThin
$ws.Range($styleF).Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range($styleF).Borders.Item($xlEdgeBottom).Weight = $xlThin

# ---------------------------------------------------------------------------
# Row 20 : "Epoch 1 Min Loss" - re-styled like the header row (D2/E2..) and
# filled in with real data.
# ---------------------------------------------------------------------------
Copy-Format "D2" "D20"
$ws.Range("D20").Value = "Epoch 1 Min Loss"

Copy-Format $styleA "E20"
$ws.Range("E20").Value = 3068.42

Copy-Format $styleB "F20"
$ws.Range("F20").Value = 3052.0412999999999
Copy-Format $styleB "G20"
Copy-Format $styleB "H20"
Copy-Format $styleB "I20"

# ---------------------------------------------------------------------------
# Rows 21-28 : new "Loss 1".."Loss 8" rows
# ---------------------------------------------------------------------------
$lossLabels = @("Loss 1", "Loss 2", "Loss 3", "Loss 4", "Loss 5", "Loss 6", "Loss 7", "Loss 8")
$eVals = @(7239.3440000000001, 6803.4939999999997, 6201.4008000000003, 6005.973, 5897.692, 5821.1379999999999, 5769.4511000000002, 5721.3071)
$fVals = @(7240.6790000000001, 6808.009, 6191.0739999999996, 5998.9139999999998, 5897.4059999999999, 5826.5320000000002, 5770.317, 5727.1704)

for ($i = 0; $i -lt 8; $i++) {
    $row = 21 + $i

    # D column: same visual style as D2 (bold header look, thin border all around)
    Copy-Format "D2" "D$row"
    $ws.Range("D$row").Value = $lossLabels[$i]

    # E column: plain / default style, no border
    Copy-Format $blank "E$row"
    $ws.Range("E$row").Value = $eVals[$i]

    # F:I columns: thin left/right only
    Copy-Format $styleC "F$row"
    $ws.Range("F$row").Value = $fVals[$i]
    Copy-Format $styleC "G$row"
    Copy-Format $styleC "H$row"
    Copy-Format $styleC "I$row"
}

# ---------------------------------------------------------------------------
# Row 29 / 30 : the old "Epoch 2 Min Loss" & "Output Grade (A-F)" rows,
# moved down; row 29 also picks up real data.
# ---------------------------------------------------------------------------

# D29 <- old D21 ("Epoch 2 Min Loss"), same look as D4 above (style 10)
Copy-Format "D4" "D29"
$ws.Range("D29").Value = "Epoch 2 Min Loss"

# D30 <- old D22 ("Output Grade (A-F)"), same look as D5 above (style 10)
Copy-Format "D5" "D30"
$ws.Range("D30").Value = "Output Grade (A-F)"

Copy-Format $styleD "E29"
$ws.Range("E29").Value = 5681.0560999999998

Copy-Format $styleE "E30"

Copy-Format $styleC "F29"
$ws.Range("F29").Value = 5679.4859999999999
Copy-Format $styleC "G29"
Copy-Format $styleC "H29"
Copy-Format $styleC "I29"

Copy-Format $styleF "F30"
Copy-Format $styleF "G30"
Copy-Format $styleF "H30"
Copy-Format $styleF "I30"

# ---------------------------------------------------------------------------
# Clean up scratch cells used to build the new border styles so they don't
# show up as stray content / expand the used range.
# ---------------------------------------------------------------------------
$ws.Range("Z100:Z106").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Selection moves to F21, matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("F21").Select() | Out-Null
